$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: update I8, add J8
$ws.Range("I8").Value = 0.468739117536998
$ws.Range("J8").Value = 0.2150495036779461

# Row 9: update H9, add I9
$ws.Range("H9").Value = 0.5099036351493167
$ws.Range("I9").Value = 0.24

# Row 10: update G10, add H10
$ws.Range("G10").Value = 0.5604363747513331
$ws.Range("H10").Value = 0.3087982760018804

# Row 11: update F11, add G11
$ws.Range("F11").Value = 0.5999036351493168
$ws.Range("G11").Value = 0.32

# Row 12: update E12, add F12
$ws.Range("E12").Value = 0.6299036351493167
$ws.Range("F12").Value = 0.4476495795507702

# Row 13: update D13, add E13
$ws.Range("D13").Value = 0.3603773643037867
$ws.Range("E13").Value = 0.1088966743764388

# Row 14: update C14, add D14
$ws.Range("C14").Value = 0.4107440146302961
$ws.Range("D14").Value = 0.1461563307127136

# Row 15: update B15, add C15
$ws.Range("B15").Value = 0.25708246933236
$ws.Range("C15").Value = 0.09547648014918764

# Row 16: add B16
$ws.Range("B16").Value = 0.0959495356205764
